# Omaha_Cal_Info_GA02HYPM_00001.xlsx update
# "Global Argentine reference designators and OOI Barcodes"
# Added missing OOI barcodes and serial numbers. Corrected engineering
# reference designators, changed CTDMO and PHSEN reference designators to
# nominal values.

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Asset_Cal_Info sheet -------------------------------------------------

# PHSEN / CTDMO engineering ref des corrected to the nominal RIM01 CTDMO value
$assetCal.Range("A41").Value = "GA02HYPM-RIM01-02-CTDMOG039"
$assetCal.Range("A42").Value = "GA02HYPM-RIM01-02-CTDMOG039"
$assetCal.Range("A43").Value = "GA02HYPM-RIM01-02-CTDMOG039"

# GP001 engineering ref des corrected to the nominal RIM01 SIO engineering ref des
$assetCal.Range("A45").Value = "GA02HYPM-RIM01-00-SIOENG000"
# Missing OOI barcode added for the RIM01 SIO engineering sensor
$assetCal.Range("E45").Value = "OL000005"

# WFP02/WFP03 engineering ref des corrected to the nominal WFPENG value
$assetCal.Range("A46").Value = "GA02HYPM-WFP02-00-WFPENG000"
$assetCal.Range("A47").Value = "GA02HYPM-WFP03-00-WFPENG000"

$wb.Save()
